# Update the "dSF" column (F) values for a set of rows.
# These correspond to a "repull data / push all data / mean calculation"
# refresh where the dSF figures were recomputed and now differ from
# the previously-pulled dS0 values in several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -1
    9  = 1
    15 = -2
    17 = -3
    19 = -4
    20 = -5
    21 = 8
    22 = -4
    24 = 1
    25 = -7
    29 = -4
    30 = 3
    34 = 1
    39 = -3
    42 = -1
    48 = -1
    50 = -1
    52 = 3
    58 = 1
    61 = -1
    62 = 5
    63 = 3
    65 = -1
    66 = -3
    71 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
